$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Re-style existing row 8 so it becomes a "group bottom" row (same look as
#    row 4): adds a bordered A8 cell and switches B8:E8 to styles 8/9.
# ---------------------------------------------------------------------------
$ws.Range("A4:E4").Copy()
$ws.Range("A8:E8").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Apply formatting (borders/fonts via style copy) to the new rows 9-16
#    BEFORE filling in values, so the new shared strings are appended to
#    sharedStrings.xml strictly in the order the values are assigned below.
# ---------------------------------------------------------------------------
$ws.Range("A4:E4").Copy()
$ws.Range("A9:E9").PasteSpecial(-4122)

$ws.Range("B5:E5").Copy()
$ws.Range("B10:E10").PasteSpecial(-4122)
$ws.Range("B11:E11").PasteSpecial(-4122)
$ws.Range("B12:E12").PasteSpecial(-4122)

$ws.Range("A3:E3").Copy()
$ws.Range("A13:E13").PasteSpecial(-4122)

$ws.Range("B5:E5").Copy()
$ws.Range("B14:E14").PasteSpecial(-4122)
$ws.Range("B15:E15").PasteSpecial(-4122)
$ws.Range("B16:E16").PasteSpecial(-4122)

# Row heights that differ from the sheet default.
$ws.Rows.Item(9).RowHeight = 21.6
$ws.Rows.Item(11).RowHeight = 26.4
$ws.Rows.Item(12).RowHeight = 21.6
$ws.Rows.Item(13).RowHeight = 31.8
$ws.Rows.Item(15).RowHeight = 21.6
$ws.Rows.Item(16).RowHeight = 21.6

# ---------------------------------------------------------------------------
# 3) Fill in the numbers (column B).
# ---------------------------------------------------------------------------
$ws.Range("B9").Value = 161
$ws.Range("B10").Value = 97
$ws.Range("B11").Value = 100
$ws.Range("B12").Value = 103
$ws.Range("B13").Value = 106
$ws.Range("B14").Value = 67
$ws.Range("B15").Value = 70
$ws.Range("B16").Value = 73

# ---------------------------------------------------------------------------
# 4) Fill in the text (columns C/D/E), in the exact order the new shared
#    strings were originally introduced.
# ---------------------------------------------------------------------------

# Row 9 : C, D, E
$ws.Range("C9").Value = " ...[K]Team [CS:X]Charm[CR]... Kwah-ha-ha...\nJust lovely…"
$ws.Range("D9").Value = " ...[K]Команда [CS:X]Шарм[CR]... Квох-ха-ха...\nПросто чудесны..."
$ws.Range("E9").Value = " ...[K]Ëïíàîäà [CS:X]Šàñí[CR]... Ëâïö-öà-öà...\nÐñïòóï œôäåòîú..."

# Column C, rows 10-13 (English)
$ws.Range("C10").Value = " ...[K]Malevolent darkness…"
$ws.Range("C11").Value = " The darkness spreads and\ndeepens…"
$ws.Range("C12").Value = " Its progress can be halted,\nyes...[K] But only by you."
$ws.Range("C13").Value = " I shall ask this of you.[K]\nStop [CS:N]Darkrai[CR]'s future of nightmares…"

# Column D, rows 10-13 (Russian)
$ws.Range("D10").Value = " ...[K]Зловещая тьма..."
$ws.Range("D11").Value = " Тьма ширится и укрепляется..."
$ws.Range("D12").Value = " Её можно остановить, да...[K]\nНо только вам это по силам."
$ws.Range("D13").Value = " Я прошу вас.[K] Не дайте сбыться\nкошмарному будущему [CS:N]Даркрая[CR]..."

# Column C, rows 14-16 (English)
$ws.Range("C14").Value = " ...[K]The darkness...`""
$ws.Range("C15").Value = " The malevolent darkness has\nmelted away…"
$ws.Range("C16").Value = " It is your doing. You have\nmy thanks."

# Column D, rows 14-16 (Russian)
$ws.Range("D14").Value = " ...[K]Тьма..."
$ws.Range("D15").Value = " Зловешая тьма отступила..."
$ws.Range("D16").Value = " Это ваша заслуга. Искренне\nблагодарю вас."

# Column E, rows 14-16 (corrupted encoding)
$ws.Range("E14").Value = " ...[K]Óûíà..."
$ws.Range("E15").Value = " Èìïâåšàÿ óûíà ïóòóôðéìà..."
$ws.Range("E16").Value = " Üóï âàšà èàòìôãà. Éòëñåîîå\náìàãïäàñý âàò."

# Column E, rows 10-13 (corrupted encoding) - filled in last.
$ws.Range("E10").Value = " ...[K]Èìïâåþàÿ óûíà..."
$ws.Range("E11").Value = " Óûíà šéñéóòÿ é ôëñåðìÿåóòÿ..."
$ws.Range("E12").Value = " Åæ íïçîï ïòóàîïâéóû, äà...[K]\nÎï óïìûëï âàí üóï ðï òéìàí."
$ws.Range("E13").Value = " Ÿ ðñïšô âàò.[K] Îå äàêóå òáúóûòÿ\nëïšíàñîïíô áôäôþåíô [CS:N]Äàñëñàÿ[CR]..."

# ---------------------------------------------------------------------------
# 5) Selection follows the edit (matches the author ending on D13).
# ---------------------------------------------------------------------------
[void]$ws.Range("D13").Select()
